# Append a new sentence (as a brand-new run) right after the existing
# run that ends with the text "...leggere il "LEGGIMI.txt"."
#
# The sentence to append is:
#   " Attenzione: modificare entrambi i parametri max_allowed_packet
#   impostandoli a 100M, soprattutto quello sotto le impostazioni del server!"
#
# We locate the unique text "LEGGIMI.txt”." in the document, collapse the
# found range to its end point, and insert the new text there. Because the
# insertion point sits right after the end of the existing run (outside of
# it), Word creates a brand-new run for the inserted text instead of
# merging it into the previous run - matching the target markup exactly.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("LEGGIMI.txt”.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Collapse(0)
$rng.InsertAfter(" Attenzione: modificare entrambi i parametri max_allowed_packet impostandoli a 100M, soprattutto quello sotto le impostazioni del server!")
